$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 467.28
$ws.Range("J17").Value = 467.45456
$ws.Range("L17").Value = 1402.36368
$ws.Range("N17").Value = -1738.36368
# Row 125
$ws.Range("H125").Value = 2050
$ws.Range("I125").Value = 1333.3334
$ws.Range("J125").Value = 4200
$ws.Range("K125").Value = 12000.0006
$ws.Range("L125").Value = 37800
$ws.Range("M125").Value = -9540.000599999999
$ws.Range("N125").Value = -42720
# Row 138
$ws.Range("H138").Value = 2334.65
$ws.Range("I138").Value = 1896.6216
$ws.Range("J138").Value = 3039.3044
$ws.Range("K138").Value = 5689.864799999999
$ws.Range("L138").Value = 9117.913199999999
$ws.Range("M138").Value = -549.8647999999994
$ws.Range("N138").Value = -19397.9132

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 818006.4399999999
$ws.Range("I2").Value = 735.56525
$ws.Range("K2").Value = 735.56525
$ws.Range("M2").Value = -622.56525
# Row 32
$ws.Range("H32").Value = 6017.8945
$ws.Range("I32").Value = 4366.2793
$ws.Range("J32").Value = 21800
$ws.Range("K32").Value = 4366.2793
$ws.Range("L32").Value = 21800
$ws.Range("M32").Value = -4079.2793
$ws.Range("N32").Value = -22374
# Row 97
$ws.Range("H97").Value = 847.64703
$ws.Range("I97").Value = 528.5
$ws.Range("J97").Value = 2337
$ws.Range("K97").Value = 528.5
$ws.Range("L97").Value = 2337
$ws.Range("M97").Value = -32.5
$ws.Range("N97").Value = -3329
# Row 101
$ws.Range("H101").Value = 36663.168
$ws.Range("J101").Value = 36663.168
$ws.Range("L101").Value = 36663.168
$ws.Range("N101").Value = -43153.168
# Row 102
$ws.Range("H102").Value = 1977.5883
$ws.Range("I102").Value = 1964.5385
$ws.Range("J102").Value = 2020
$ws.Range("K102").Value = 1964.5385
$ws.Range("L102").Value = 2020
$ws.Range("M102").Value = -342.5385000000001
$ws.Range("N102").Value = -5264
# Row 110
$ws.Range("H110").Value = 1642.2222
$ws.Range("I110").Value = 1285.7142
$ws.Range("K110").Value = 1285.7142
$ws.Range("M110").Value = 759.2858000000001
# Row 116
$ws.Range("H116").Value = 818006.4399999999
$ws.Range("I116").Value = 735.56525
$ws.Range("K116").Value = 735.56525
$ws.Range("M116").Value = 1558.43475
# Row 132
$ws.Range("H132").Value = 711270
$ws.Range("I132").Value = 842707.0600000001
$ws.Range("J132").Value = 112501
$ws.Range("K132").Value = 2528121.18
$ws.Range("L132").Value = 337503
$ws.Range("M132").Value = -2525591.18
$ws.Range("N132").Value = -342563

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 818006.4399999999
$ws.Range("I3").Value = 735.56525
$ws.Range("K3").Value = 735.56525
$ws.Range("M3").Value = -621.56525
# Row 86
$ws.Range("H86").Value = 2168.6
$ws.Range("I86").Value = 1788.8695
$ws.Range("J86").Value = 2682.353
$ws.Range("K86").Value = 1788.8695
$ws.Range("L86").Value = 2682.353
$ws.Range("M86").Value = -665.8695
$ws.Range("N86").Value = -4928.353
# Row 89
$ws.Range("H89").Value = 2168.6
$ws.Range("I89").Value = 1788.8695
$ws.Range("J89").Value = 2682.353
$ws.Range("K89").Value = 8944.3475
$ws.Range("L89").Value = 13411.765
$ws.Range("M89").Value = -3328.3475
$ws.Range("N89").Value = -24643.765
# Row 105
$ws.Range("H105").Value = 1470.5
$ws.Range("I105").Value = 1428.9412
$ws.Range("J105").Value = 1571.4286
$ws.Range("K105").Value = 1428.9412
$ws.Range("L105").Value = 1571.4286
$ws.Range("M105").Value = 318.0588
$ws.Range("N105").Value = -5065.4286
# Row 134
$ws.Range("H134").Value = 10585139
$ws.Range("I134").Value = 11830244
$ws.Range("J134").Value = 1750
$ws.Range("K134").Value = 35490732
$ws.Range("L134").Value = 5250
$ws.Range("M134").Value = -35488197
$ws.Range("N134").Value = -10320

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 359560.1
$ws.Range("I31").Value = 1169.0667
$ws.Range("K31").Value = 1169.0667
$ws.Range("M31").Value = -874.0667000000001
# Row 34
$ws.Range("H34").Value = 359560.1
$ws.Range("I34").Value = 1169.0667
$ws.Range("K34").Value = 1169.0667
$ws.Range("M34").Value = -967.0667000000001
# Row 105
$ws.Range("H105").Value = 902.8889
$ws.Range("I105").Value = 804.2857
$ws.Range("J105").Value = 1248
$ws.Range("K105").Value = 804.2857
$ws.Range("L105").Value = 1248
$ws.Range("M105").Value = 942.7143
$ws.Range("N105").Value = -4742
# Row 132
$ws.Range("H132").Value = 2238.0293
$ws.Range("I132").Value = 2102.5483
$ws.Range("J132").Value = 3638
$ws.Range("K132").Value = 6307.644899999999
$ws.Range("L132").Value = 10914
$ws.Range("M132").Value = -3777.644899999999
$ws.Range("N132").Value = -15974

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 57
$ws.Range("H57").Value = 2000
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
# Row 113
$ws.Range("H113").Value = 449.98935
$ws.Range("I113").Value = 436.85245
$ws.Range("K113").Value = 1310.55735
$ws.Range("M113").Value = 859.44265
# Row 122
$ws.Range("H122").Value = 9921711
$ws.Range("I122").Value = 16667215
$ws.Range("J122").Value = 1489831.2
$ws.Range("K122").Value = 150004935
$ws.Range("L122").Value = 13408480.8
$ws.Range("M122").Value = -150002485
$ws.Range("N122").Value = -13413380.8
# Row 131
$ws.Range("H131").Value = 3332.4546
$ws.Range("I131").Value = 6507
$ws.Range("J131").Value = 2398.7646
$ws.Range("K131").Value = 19521
$ws.Range("L131").Value = 7196.293799999999
$ws.Range("M131").Value = -14481
$ws.Range("N131").Value = -17276.2938

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 4878.5557
$ws.Range("I70").Value = 4817.8335
$ws.Range("J70").Value = 5000
$ws.Range("K70").Value = 4817.8335
$ws.Range("L70").Value = 5000
$ws.Range("M70").Value = -4547.8335
$ws.Range("N70").Value = -5540
# Row 73
$ws.Range("H73").Value = 4878.5557
$ws.Range("I73").Value = 4817.8335
$ws.Range("J73").Value = 5000
$ws.Range("K73").Value = 4817.8335
$ws.Range("L73").Value = 5000
$ws.Range("M73").Value = -3881.8335
$ws.Range("N73").Value = -6872
# Row 97
$ws.Range("H97").Value = 1657.3334
$ws.Range("I97").Value = 1638.3334
$ws.Range("J97").Value = 1733.3334
$ws.Range("K97").Value = 1638.3334
$ws.Range("L97").Value = 1733.3334
$ws.Range("M97").Value = -1142.3334
$ws.Range("N97").Value = -2725.3334
# Row 102
$ws.Range("H102").Value = 1604.8918
$ws.Range("I102").Value = 1504.6786
$ws.Range("J102").Value = 1916.6666
$ws.Range("K102").Value = 1504.6786
$ws.Range("L102").Value = 1916.6666
$ws.Range("M102").Value = 117.3214
$ws.Range("N102").Value = -5160.6666
# Row 126
$ws.Range("H126").Value = 2304.5
$ws.Range("I126").Value = 2130.75
$ws.Range("J126").Value = 2999.5
$ws.Range("K126").Value = 6392.25
$ws.Range("L126").Value = 8998.5
$ws.Range("M126").Value = -3922.25
$ws.Range("N126").Value = -13938.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 3326.1904
$ws.Range("I16").Value = 2671.6667
$ws.Range("J16").Value = 7253.3335
$ws.Range("K16").Value = 2671.6667
$ws.Range("L16").Value = 7253.3335
$ws.Range("M16").Value = -2501.6667
$ws.Range("N16").Value = -7593.3335
# Row 40
$ws.Range("H40").Value = 2575
$ws.Range("I40").Value = 2266.6667
$ws.Range("J40").Value = 3500
$ws.Range("K40").Value = 2266.6667
$ws.Range("L40").Value = 3500
$ws.Range("M40").Value = -2130.6667
$ws.Range("N40").Value = -3772

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 4750.143
$ws.Range("I81").Value = 1398.8572
$ws.Range("J81").Value = 8101.4287
$ws.Range("K81").Value = 2797.7144
$ws.Range("L81").Value = 16202.8574
$ws.Range("M81").Value = -1736.7144
$ws.Range("N81").Value = -18324.8574
# Row 84
$ws.Range("H84").Value = 4750.143
$ws.Range("I84").Value = 1398.8572
$ws.Range("J84").Value = 8101.4287
$ws.Range("K84").Value = 13988.572
$ws.Range("L84").Value = 81014.28700000001
$ws.Range("M84").Value = -8684.572
$ws.Range("N84").Value = -91622.28700000001
# Row 100
$ws.Range("H100").Value = 527.1429000000001
$ws.Range("I100").Value = 458
$ws.Range("J100").Value = 700
$ws.Range("K100").Value = 916
$ws.Range("L100").Value = 1400
$ws.Range("M100").Value = -375
$ws.Range("N100").Value = -2482
# Row 132
$ws.Range("H132").Value = 3873.4468
$ws.Range("I132").Value = 4296.951
$ws.Range("K132").Value = 12890.853
$ws.Range("M132").Value = -10360.853
